$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newA = "(003)_b2c_dsl_vasaktivate_aktivuj dsl + rental + security na b2c"
$newB = "vasaktivate"

for ($r = 2; $r -le 24; $r++) {
    $ws.Cells.Item($r, 1).Value = $newA
    $ws.Cells.Item($r, 2).Value = $newB
}
